$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.093.01"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.795.05"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Formula = "'601.40"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Formula = "'165.05"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Formula = "'6.49"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Formula = "'35.80"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "4.431.82"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "3.790.46"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "68.082.84"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Formula = "'18.38"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Formula = "'7.08"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Formula = "'461.39"
$ws.Range("D21").Formula = "'9.71"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -4.76%  "
$ws.Range("D24").Formula = "'82.99"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Formula = "'12.01"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "3.943.54"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").Formula = "'7.35"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Formula = "'29.33"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Formula = "'9.03"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Formula = "'0.0996"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Formula = "'5.84"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Formula = "'1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Formula = "'47.57"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Formula = "'0.300"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Formula = "'43.29"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Formula = "'151.91"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Formula = "'391.82"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Formula = "'26.56"
$ws.Range("E51").Value = "  -1.01%  "